$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 32; existing rows 32-38 shift down to 33-39.
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with a new weekly price record
# (same market/category/quality metadata as the surrounding rows, new date
# and the price figures that used to belong to the old row 32 shifted down).
$ws.Range("A32").Value = 8
$ws.Range("B32").Value = "Terminal La Palmera de La Serena"
$ws.Range("C32").Value = "Coquimbo"
$ws.Range("D32").Value = 44463
$ws.Range("E32").Value = 4
$ws.Range("F32").Value = 100112052
$ws.Range("G32").Value = "Albahaca"
$ws.Range("H32").Value = "Sin especificar"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 800
$ws.Range("K32").Value = 4000
$ws.Range("L32").Value = 4500
$ws.Range("M32").Value = 4250
$ws.Range("N32").Value = "$/paquete"
$ws.Range("O32").Value = "Región de Arica y Parinacota"
$ws.Range("P32").Value = 4250
$ws.Range("Q32").Value = 1
$ws.Range("R32").Value = "Hortaliza"
